$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1983.3334
$ws.Range("I40").Value = 1983.3334
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1983.3334
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1808.3334
$ws.Range("N40").ClearContents()
$ws.Range("H51").Value = 1650
$ws.Range("I51").Value = 1650
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1650
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1166
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 7084.933
$ws.Range("I62").Value = 6104.1665
$ws.Range("J62").Value = 7738.778
$ws.Range("K62").Value = 6104.1665
$ws.Range("L62").Value = 7738.778
$ws.Range("M62").Value = -5480.1665
$ws.Range("N62").Value = -8986.778
$ws.Range("H64").Value = 3266.6667
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 3100
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 3100
$ws.Range("M64").Value = -3352
$ws.Range("N64").Value = -3596
$ws.Range("H65").Value = 7084.933
$ws.Range("I65").Value = 6104.1665
$ws.Range("J65").Value = 7738.778
$ws.Range("K65").Value = 30520.8325
$ws.Range("L65").Value = 38693.89
$ws.Range("M65").Value = -27400.8325
$ws.Range("N65").Value = -44933.89
$ws.Range("H67").Value = 3266.6667
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 3100
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 3100
$ws.Range("M67").Value = -2742
$ws.Range("N67").Value = -4816
$ws.Range("H80").Value = 310.8125
$ws.Range("I80").Value = 307.5
$ws.Range("J80").Value = 320.75
$ws.Range("K80").Value = 922.5
$ws.Range("L80").Value = 962.25
$ws.Range("M80").Value = 75.5
$ws.Range("N80").Value = -2958.25
$ws.Range("H83").Value = 310.8125
$ws.Range("I83").Value = 307.5
$ws.Range("J83").Value = 320.75
$ws.Range("K83").Value = 2767.5
$ws.Range("L83").Value = 2886.75
$ws.Range("M83").Value = 2224.5
$ws.Range("N83").Value = -12870.75
$ws.Range("H112").Value = 2279.8125
$ws.Range("I112").Value = 1805
$ws.Range("J112").Value = 2389.3845
$ws.Range("K112").Value = 5415
$ws.Range("L112").Value = 7168.1535
$ws.Range("M112").Value = -4307
$ws.Range("N112").Value = -9384.1535
$ws.Range("H113").Value = 5595.75
$ws.Range("I113").Value = 4966.857
$ws.Range("J113").Value = 9998
$ws.Range("K113").Value = 4966.857
$ws.Range("L113").Value = 9998
$ws.Range("M113").Value = -1712.857
$ws.Range("N113").Value = -16506
$ws.Range("H137").Value = 6017.4375
$ws.Range("I137").Value = 6798.7
$ws.Range("J137").Value = 4715.3335
$ws.Range("K137").Value = 20396.1
$ws.Range("L137").Value = 14146.0005
$ws.Range("M137").Value = -17846.1
$ws.Range("N137").Value = -19246.0005
$ws.Range("H138").Value = 9072.682000000001
$ws.Range("I138").Value = 10944.467
$ws.Range("J138").Value = 8552.741
$ws.Range("K138").Value = 32833.401
$ws.Range("L138").Value = 25658.223
$ws.Range("M138").Value = -27693.401
$ws.Range("N138").Value = -35938.223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16764.527
$ws.Range("I32").Value = 10621.72
$ws.Range("J32").Value = 30725.455
$ws.Range("K32").Value = 10621.72
$ws.Range("L32").Value = 30725.455
$ws.Range("M32").Value = -10334.72
$ws.Range("N32").Value = -31299.455
$ws.Range("H102").Value = 2239.2
$ws.Range("I102").Value = 2168.2
$ws.Range("J102").Value = 2274.7
$ws.Range("K102").Value = 2168.2
$ws.Range("L102").Value = 2274.7
$ws.Range("M102").Value = -546.1999999999998
$ws.Range("N102").Value = -5518.7
$ws.Range("H112").Value = 54021.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 54021.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 54021.5
$ws.Range("N112").Value = -56975.5
$ws.Range("H135").Value = 249998
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 249998
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 249998
$ws.Range("N135").Value = -260138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H54").Value = 4272.5
$ws.Range("I54").Value = 4272.5
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 4272.5
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -3788.5
$ws.Range("H86").Value = 7141
$ws.Range("I86").Value = 7201.25
$ws.Range("J86").Value = 6900
$ws.Range("K86").Value = 7201.25
$ws.Range("L86").Value = 6900
$ws.Range("M86").Value = -6078.25
$ws.Range("N86").Value = -9146
$ws.Range("H89").Value = 7141
$ws.Range("I89").Value = 7201.25
$ws.Range("J89").Value = 6900
$ws.Range("K89").Value = 36006.25
$ws.Range("L89").Value = 34500
$ws.Range("M89").Value = -30390.25
$ws.Range("N89").Value = -45732
$ws.Range("H99").Value = 988.0833
$ws.Range("I99").Value = 998.1579
$ws.Range("J99").Value = 949.8
$ws.Range("K99").Value = 998.1579
$ws.Range("L99").Value = 949.8
$ws.Range("M99").Value = 499.8421
$ws.Range("N99").Value = -3945.8
$ws.Range("H135").Value = 500000000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 500000000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 500000000
$ws.Range("N135").Value = -500010140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 30000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30782
$ws.Range("M39").ClearContents()
$ws.Range("H49").Value = 30000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30364
$ws.Range("M49").ClearContents()
$ws.Range("H51").Value = 59999
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 59999
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 59999
$ws.Range("N51").Value = -61471
$ws.Range("H61").Value = 59999
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 59999
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 59999
$ws.Range("N61").Value = -60695
$ws.Range("H132").Value = 2268.2856
$ws.Range("I132").Value = 2094.3215
$ws.Range("J132").Value = 2964.1428
$ws.Range("K132").Value = 6282.9645
$ws.Range("L132").Value = 8892.428400000001
$ws.Range("M132").Value = -3752.9645
$ws.Range("N132").Value = -13952.4284
$ws.Range("H141").Value = 52698.855
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 52698.855
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 52698.855
$ws.Range("N141").Value = -63058.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 847.0769
$ws.Range("I33").Value = 2060.4
$ws.Range("J33").Value = 88.75
$ws.Range("K33").Value = 12362.4
$ws.Range("L33").Value = 532.5
$ws.Range("M33").Value = -12079.4
$ws.Range("N33").Value = -1098.5
$ws.Range("H113").Value = 4313.636
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 4420
$ws.Range("K113").Value = 9750
$ws.Range("L113").Value = 13260
$ws.Range("M113").Value = -7580
$ws.Range("N113").Value = -17600
$ws.Range("H140").Value = 3841.9333
$ws.Range("I140").Value = 3356.077
$ws.Range("J140").Value = 7000
$ws.Range("K140").Value = 10068.231
$ws.Range("L140").Value = 21000
$ws.Range("M140").Value = -4888.231
$ws.Range("N140").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1512.25
$ws.Range("I3").Value = 275
$ws.Range("J3").Value = 2749.5
$ws.Range("K3").Value = 275
$ws.Range("L3").Value = 2749.5
$ws.Range("M3").Value = -159
$ws.Range("N3").Value = -2981.5
$ws.Range("H11").Value = 3275247.2
$ws.Range("I11").Value = 12001667
$ws.Range("J11").Value = 2840
$ws.Range("K11").Value = 12001667
$ws.Range("L11").Value = 2840
$ws.Range("M11").Value = -12001528
$ws.Range("N11").Value = -3118
$ws.Range("H24").Value = 21620.777
$ws.Range("I24").Value = 18555
$ws.Range("J24").Value = 22004
$ws.Range("K24").Value = 18555
$ws.Range("L24").Value = 22004
$ws.Range("M24").Value = -18382
$ws.Range("N24").Value = -22350
$ws.Range("H29").Value = 17999.928
$ws.Range("I29").Value = 16200
$ws.Range("J29").Value = 18999.889
$ws.Range("K29").Value = 16200
$ws.Range("L29").Value = 18999.889
$ws.Range("M29").Value = -15910
$ws.Range("N29").Value = -19579.889
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4863.375
$ws.Range("I7").Value = 4801.8
$ws.Range("J7").Value = 4966
$ws.Range("K7").Value = 4801.8
$ws.Range("L7").Value = 4966
$ws.Range("M7").Value = -4689.8
$ws.Range("N7").Value = -5190
$ws.Range("H16").Value = 3357.9333
$ws.Range("I16").Value = 3582.2307
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 3582.2307
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -3412.2307
$ws.Range("N16").Value = -2240
$ws.Range("H40").Value = 4333.3335
$ws.Range("I40").Value = 4333.3335
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4333.3335
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4197.3335
$ws.Range("H46").Value = 2941.1177
$ws.Range("I46").Value = 2800
$ws.Range("J46").Value = 2999.9167
$ws.Range("K46").Value = 2800
$ws.Range("L46").Value = 2999.9167
$ws.Range("M46").Value = -2612
$ws.Range("N46").Value = -3375.9167
$ws.Range("H110").Value = 60214.668
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 60214.668
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 60214.668
$ws.Range("N110").Value = -68394.66800000001
$ws.Range("H126").Value = 4863.375
$ws.Range("I126").Value = 4801.8
$ws.Range("J126").Value = 4966
$ws.Range("K126").Value = 14405.4
$ws.Range("L126").Value = 14898
$ws.Range("M126").Value = -11935.4
$ws.Range("N126").Value = -19838
$ws.Range("H136").Value = 2671.2856
$ws.Range("I136").Value = 2816.6667
$ws.Range("J136").Value = 1799
$ws.Range("K136").Value = 8450.000100000001
$ws.Range("L136").Value = 5397
$ws.Range("M136").Value = -5900.000100000001
$ws.Range("N136").Value = -10497

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 69436.266
$ws.Range("I14").Value = 127200.5
$ws.Range("J14").Value = 3420
$ws.Range("K14").Value = 127200.5
$ws.Range("L14").Value = 3420
$ws.Range("M14").Value = -127032.5
$ws.Range("N14").Value = -3756
$ws.Range("H81").Value = 3849.5715
$ws.Range("I81").Value = 3607.2307
$ws.Range("J81").Value = 7000
$ws.Range("K81").Value = 7214.4614
$ws.Range("L81").Value = 14000
$ws.Range("M81").Value = -6153.4614
$ws.Range("N81").Value = -16122
$ws.Range("H84").Value = 3849.5715
$ws.Range("I84").Value = 3607.2307
$ws.Range("J84").Value = 7000
$ws.Range("K84").Value = 36072.307
$ws.Range("L84").Value = 70000
$ws.Range("M84").Value = -30768.307
$ws.Range("N84").Value = -80608
$ws.Range("H132").Value = 2318.75
$ws.Range("I132").Value = 966
$ws.Range("J132").Value = 3285
$ws.Range("K132").Value = 2898
$ws.Range("L132").Value = 9855
$ws.Range("M132").Value = -368
$ws.Range("N132").Value = -14915

